$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 14:05"

# Country-name cell swaps caused by reordering in the shared-string table
$ws.Range("A51").Value = "Barein"
$ws.Range("A52").Value = "Dinamarca"
$ws.Range("A53").Value = "Corea del Sur"
$ws.Range("A54").Value = "Serbia"
$ws.Range("A81").Value = "Consejo Danes para los Refugiados"
$ws.Range("A82").Value = "Tailandia"
$ws.Range("A89").Value = "Republica de Macedonia"
$ws.Range("A90").Value = "Croacia"
$ws.Range("A109").Value = "Etiopia"
$ws.Range("A110").Value = "Guinea-Bisau"
$ws.Range("A111").Value = "Libano"

# Updated numeric statistics per row
$ws.Range("B4").Value = 1837625
$ws.Range("C4").Value = 455
$ws.Range("E4").Value = 1131552
$ws.Range("B10").Value = 191333
$ws.Range("C10").Value = 724
$ws.Range("E10").Value = 93893
$ws.Range("B23").Value = 58433
$ws.Range("C23").Value = 1523
$ws.Range("D23").Value = 33437
$ws.Range("E23").Value = 24956
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 40
$ws.Range("B25").Value = 46545
$ws.Range("C25").Value = 103
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 5962
$ws.Range("B26").Value = 43403
$ws.Range("C26").Value = 847
$ws.Range("D26").Value = 18776
$ws.Range("E26").Value = 24387
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 240
$ws.Range("B28").Value = 37814
$ws.Range("C28").Value = 272
$ws.Range("E28").Value = 28440
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 4403
$ws.Range("B47").Value = 16733
$ws.Range("C47").Value = 2
$ws.Range("D47").Value = 15596
$ws.Range("E47").Value = 469
$ws.Range("B51").Value = 11804
$ws.Range("C51").Value = 406
$ws.Range("D51").Value = 7070
$ws.Range("E51").Value = 4715
$ws.Range("H51").Value = 19
$ws.Range("B52").Value = 11699
$ws.Range("C52").Value = 30
$ws.Range("D52").Value = 10412
$ws.Range("E52").Value = 711
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 576
$ws.Range("B53").Value = 11503
$ws.Range("C53").Value = 35
$ws.Range("D53").Value = 10422
$ws.Range("E53").Value = 810
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 271
$ws.Range("B54").Value = 11412
$ws.Range("D54").Value = 6698
$ws.Range("E54").Value = 4471
$ws.Range("H54").Value = 243
$ws.Range("D55").Value = 5458
$ws.Range("E55").Value = 5809
$ws.Range("D79").Value = 2859
$ws.Range("E79").Value = 788
$ws.Range("B81").Value = 3195
$ws.Range("C81").Value = 125
$ws.Range("D81").Value = 454
$ws.Range("E81").Value = 2669
$ws.Range("H81").Value = 72
$ws.Range("B82").Value = 3082
$ws.Range("C82").Value = 1
$ws.Range("D82").Value = 2965
$ws.Range("E82").Value = 60
$ws.Range("H82").Value = 57
$ws.Range("B89").Value = 2315
$ws.Range("C89").Value = 89
$ws.Range("D89").Value = 1569
$ws.Range("E89").Value = 606
$ws.Range("G89").Value = 7
$ws.Range("H89").Value = 140
$ws.Range("B90").Value = 2246
$ws.Range("D90").Value = 2077
$ws.Range("E90").Value = 66
$ws.Range("H90").Value = 103
$ws.Range("E99").Value = 1314
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 6
$ws.Range("B102").Value = 1634
$ws.Range("C102").Value = 1
$ws.Range("E102").Value = 812
$ws.Range("B109").Value = 1257
$ws.Range("C109").Value = 85
$ws.Range("D109").Value = 217
$ws.Range("E109").Value = 1028
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 12
$ws.Range("B110").Value = 1256
$ws.Range("D110").Value = 42
$ws.Range("E110").Value = 1206
$ws.Range("H110").Value = 8
$ws.Range("B111").Value = 1220
$ws.Range("D111").Value = 712
$ws.Range("E111").Value = 481
$ws.Range("H111").Value = 27
